# Updated CVDs for the month
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Juarez Holdings Ss Mexico  (style 13)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Juarez Holdings Ss Mexico")
$ws.Range("E2").Value = 0.0463
$ws.Range("E3").Value = 0.0463
$ws.Range("E4").Value = 0.0463
$ws.Range("O4:W4").Value = 0

# ---------------------------------------------------------------------------
# Sheet: Manila Philippines  (style 15)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Manila Philippines")
$ws.Range("E2").Value = 0.0497
$ws.Range("E3").Value = 0.0497
$ws.Range("E4").Value = 0.0497
$ws.Range("J4").Value = 0.0154
$ws.Range("K4").Value = 0.0234
$ws.Range("N4").Value = 0.023
$ws.Range("O4").Value = 0.0112
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("E5").Value = 0.5
$ws.Range("E6").Value = 0.5
$ws.Range("E7").Value = 0.5
$ws.Range("O7").Value = 1
$ws.Range("P7:W7").Value = 0.5

# ---------------------------------------------------------------------------
# Sheet: Milwaukee Pmc Hq Wisconsin  (style 16)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("E2").Value = 0.0789
$ws.Range("E3").Value = 0.0789
$ws.Range("E4").Value = 0.0789
$ws.Range("K4").Value = 0.0541
$ws.Range("M4").Value = 0.0263
$ws.Range("N4").Value = 0.0796
$ws.Range("O4:W4").Value = 0
$ws.Range("E5").Value = 0.5
$ws.Range("E6").Value = 0.5
$ws.Range("E7").Value = 0.5
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 0.5
$ws.Range("O7").Value = 1
$ws.Range("P7:W7").Value = 0.5
$ws.Range("E8").Value = 0.2273
$ws.Range("E9").Value = 0.2273
$ws.Range("E10").Value = 0.2273
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0.032475
$ws.Range("Q10").Value = 0.032475
$ws.Range("R10").Value = 0.097425
$ws.Range("S10").Value = 0.032475
$ws.Range("T10").Value = 0.032475
$ws.Range("U10").Value = 0.032475
$ws.Range("V10").Value = 0.097425
$ws.Range("W10").Value = 0.3897

# ---------------------------------------------------------------------------
# Sheet: Pune India  (style 21)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Pune India")
$ws.Range("E2").Value = 0.2326
$ws.Range("E3").Value = 0.2326
$ws.Range("E4").Value = 0.2326
$ws.Range("O4:W4").Value = 0

# ---------------------------------------------------------------------------
# Sheet: Ratingen Germany  (style 22)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ratingen Germany")
$ws.Range("E2").Value = 0.4762
$ws.Range("E3").Value = 0.4762
$ws.Range("E4").Value = 0.4762
$ws.Range("O4:W4").Value = 0

# ---------------------------------------------------------------------------
# Sheet: Rosemont Illinois  (style 23)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Rosemont Illinois")
$ws.Range("E2").Value = 0.0943
$ws.Range("E3").Value = 0.0943
$ws.Range("E4").Value = 0.0943
$ws.Range("H4").Value = 0.025
$ws.Range("J4").Value = 0.0248
$ws.Range("K4").Value = 0.0244
$ws.Range("M4").Value = 0.0217
$ws.Range("N4").Value = 0.0458
$ws.Range("O4").Value = 0.0222
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("E5").Value = 0.125
$ws.Range("E6").Value = 0.125
$ws.Range("E7").Value = 0.125
$ws.Range("M7").Value = 0.3333
$ws.Range("N7").Value = 0.1667
$ws.Range("O7").Value = ""
$ws.Range("P7:W7").Value = 0.125

# ---------------------------------------------------------------------------
# Sheet: Tipp City Ohio  (style 29)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tipp City Ohio")
$ws.Range("O5").Value = ""

# ---------------------------------------------------------------------------
# Sheet: Apodaca Pmc Plant 1 Mexico  (style 30)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Apodaca Pmc Plant 1 Mexico")
$ws.Range("O3").Value = ""

# ---------------------------------------------------------------------------
# Sheet: Braintree Massachusetts  (style 4)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Braintree Massachusetts")
$ws.Range("E2").Value = 0.7317
$ws.Range("E3").Value = 0.7317
$ws.Range("E4").Value = 0.7317
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1
